# Update portfolio_values_app worksheet with corrected/extended history.
#
# Summary of the change:
#   1. Rows 161-206: the "buying_power" column (C) was miscalculated and is
#      corrected to equal the "portfolio_value" column (D), which was already
#      right.
#   2. A row for 2026-01-20 (serial 46042) that had been missing is inserted
#      right before the existing 2026-01-21 (serial 46043) row, pushing the
#      old rows 207-209 down to 208-210. Those shifted rows have the same
#      buying_power/portfolio_value correction applied (C = D).
#   3. Five new trading days (serials 46048-46052) are appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix buying_power (C) to match portfolio_value (D) for rows 161-206 ---
for ($r = 161; $r -le 206; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 4).Value2
}

# --- 2. Insert the missing 2026-01-20 row before row 207, shifting 207-209 to 208-210 ---
$ws.Rows.Item(207).Insert()

# Match the date-column formatting used by the surrounding rows: copy just the
# single formatted cell (A206) onto the new A207 so we reuse the existing style
# instead of synthesising a new one.
$ws.Cells.Item(206, 1).Copy()
$ws.Cells.Item(207, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(207, 1).Value = 46042
$ws.Cells.Item(207, 2).Value = 2110.53
$ws.Cells.Item(207, 3).Value = 109383.3672497559
$ws.Cells.Item(207, 4).Value = 109383.3672497559
$ws.Cells.Item(207, 5).Value = 6796.85986328125
$ws.Cells.Item(207, 6).Value = 1.16391396522522

# --- 3. Fix buying_power (C) to match portfolio_value (D) for the shifted rows 208-210 ---
for ($r = 208; $r -le 210; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 4).Value2
}

# --- 4. Append five new trading days after the previous last row (now row 210) ---
$newRows = @(
    @(46048, 2110.53, 110817.78, 110817.78, 6950.22998046875, 1.185845732688904),
    @(46049, 2110.53, 111839.14, 111839.14, 6978.60009765625, 1.187648415565491),
    @(46050, 2110.53, 113228.49, 113228.49, 6978.02978515625, 1.201764225959778),
    @(46051, 2110.53, 113611.5,  113611.5,  6969.009765625,   1.197776913642883),
    @(46052, 2110.53, 112642.71, 112642.71, 6939.02978515625, 1.185395956039429)
)

$ws.Cells.Item(206, 1).Copy()

$destRow = 211
foreach ($rowData in $newRows) {
    $ws.Cells.Item($destRow, 1).Value = $rowData[0]
    $ws.Cells.Item($destRow, 2).Value = $rowData[1]
    $ws.Cells.Item($destRow, 3).Value = $rowData[2]
    $ws.Cells.Item($destRow, 4).Value = $rowData[3]
    $ws.Cells.Item($destRow, 5).Value = $rowData[4]
    $ws.Cells.Item($destRow, 6).Value = $rowData[5]

    # New rows below the existing table don't inherit the date-column style
    # automatically (no shift happened), so copy it explicitly (re-using the
    # single-cell clipboard contents copied above keeps style reuse clean).
    $ws.Cells.Item($destRow, 1).PasteSpecial(-4122)  # xlPasteFormats

    $destRow++
}
$excel.CutCopyMode = $false
